$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.640.26"
$ws.Range("E2").Value = '  +1.68%  '
$ws.Range("D3").Value = "'1.886.07"
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = "'247.19"
$ws.Range("E5").Value = '  +5.91%  '
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = "'0.4769"
$ws.Range("E7").Value = '  +1.92%  '
$ws.Range("D8").Value = "'0.2936"
$ws.Range("E8").Value = '  +3.75%  '
$ws.Range("D9").Value = "'0.06547"
$ws.Range("E9").Value = '  +1.35%  '
$ws.Range("D10").Value = "'22.08"
$ws.Range("E10").Value = '  +5.14%  '
$ws.Range("D11").Value = "'0.07732"
$ws.Range("D12").Value = "'97.92"
$ws.Range("E12").Value = '  +4.77%  '
$ws.Range("D13").Value = "'0.7453"
$ws.Range("E13").Value = '  +9.83%  '
$ws.Range("D14").Value = "'1.884.60"
$ws.Range("E14").Value = '  +1.44%  '
$ws.Range("D15").Value = "'5.174"
$ws.Range("E15").Value = '  +2.56%  '
$ws.Range("D16").Value = "'277.16"
$ws.Range("E16").Value = '  +4.14%  '
$ws.Range("D17").Value = "'30.592.15"
$ws.Range("E17").Value = '  +1.62%  '
$ws.Range("E18").Value = '  +2.22%  '
$ws.Range("D19").Value = "'0.000007610"
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").Value = "'2.122.02"
$ws.Range("E21").Value = '  +0.57%  '
$ws.Range("D22").Value = "'5.296"
$ws.Range("E22").Value = '  +2.84%  '
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = "'6.227"
$ws.Range("E24").Value = '  +2.15%  '
$ws.Range("D25").Value = "'9.356"
$ws.Range("E25").Value = '  +0.78%  '
$ws.Range("D26").Value = "'163.72"
$ws.Range("E26").Value = '  -0.95%  '
$ws.Range("D27").Value = "'19.02"
$ws.Range("E27").Value = '  +2.75%  '
$ws.Range("D28").Value = "'1.957"
$ws.Range("E28").Value = '  +3.93%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").Value = "'0.1003"
$ws.Range("E29").Value = '  +1.99%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = "'1.374"
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("E31").Value = '  +4.69%  '
$ws.Range("D32").Value = "'4.340"
$ws.Range("E32").Value = '  +3.36%  '
$ws.Range("D33").Value = "'4.092"
$ws.Range("E33").Value = '  +2.99%  '
$ws.Range("D34").Value = "'0.04820"
$ws.Range("E34").Value = '  +3.67%  '
$ws.Range("D35").Value = "'1.133"
$ws.Range("E35").Value = '  +1.84%  '
$ws.Range("D36").Value = "'0.7054"
$ws.Range("E36").Value = '  +2.90%  '
$ws.Range("D37").Value = "'2.717"
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("D38").Value = "'0.01879"
$ws.Range("E38").Value = '  +3.66%  '
$ws.Range("D39").Value = "'2.738"
$ws.Range("D40").Value = "'6.362"
$ws.Range("E40").Value = '  +1.62%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = "'1.979"
$ws.Range("E41").Value = '  +5.67%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = "'71.63"
$ws.Range("E42").Value = '  +1.41%  '
$ws.Range("D43").Value = "'0.4246"
$ws.Range("E43").Value = '  +5.03%  '
$ws.Range("D44").Value = "'0.8417"
$ws.Range("E44").Value = '  +1.25%  '
$ws.Range("D45").Value = "'0.9998"
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").Value = "'102.94"
$ws.Range("E46").Value = '  +1.14%  '
$ws.Range("D47").Value = "'9.342"
$ws.Range("E47").Value = '  +2.13%  '
$ws.Range("D48").Value = "'7.133"
$ws.Range("E48").Value = '  +2.98%  '
$ws.Range("D49").Value = "'35.70"
$ws.Range("E49").Value = '  +4.82%  '
$ws.Range("D50").Value = "'923.05"
$ws.Range("E50").Value = '  +0.18%  '
$ws.Range("D51").Value = "'0.3919"
$ws.Range("E51").Value = '  +5.06%  '
